$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 453
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value()
    $cell.Value = $v.AddDays(1)
}
